$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 81
$ws.Range("J2").Value = 360
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 97
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 57
$ws.Range("O2").Value = 0
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 44
$ws.Range("T2").Value = 68
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 536
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 501
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 5
$ws.Range("AA2").Value = 3
